$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6996.3687
$ws.Range("I40").Value = 34566.668
$ws.Range("J40").Value = 4633.2
$ws.Range("K40").Value = 34566.668
$ws.Range("L40").Value = 4633.2
$ws.Range("M40").Value = -34391.668
$ws.Range("N40").Value = -4983.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2828749.8
$ws.Range("I2").Value = 3770928.2
$ws.Range("K2").Value = 3770928.2
$ws.Range("M2").Value = -3770815.2
$ws.Range("H32").Value = 6024.413
$ws.Range("I32").Value = 3682.4177
$ws.Range("J32").Value = 20256.54
$ws.Range("K32").Value = 3682.4177
$ws.Range("L32").Value = 20256.54
$ws.Range("M32").Value = -3395.4177
$ws.Range("N32").Value = -20830.54
$ws.Range("H45").Value = 8995378
$ws.Range("I45").Value = 17983260
$ws.Range("K45").Value = 17983260
$ws.Range("M45").Value = -17982883
$ws.Range("H61").Value = 7349.9287
$ws.Range("I61").Value = 7531.077
$ws.Range("K61").Value = 7531.077
$ws.Range("M61").Value = -7319.077
$ws.Range("H63").Value = 4503.4585
$ws.Range("I63").Value = 2362.6428
$ws.Range("J63").Value = 7500.6
$ws.Range("K63").Value = 2362.6428
$ws.Range("L63").Value = 7500.6
$ws.Range("M63").Value = -1676.6428
$ws.Range("N63").Value = -8872.6
$ws.Range("H66").Value = 4503.4585
$ws.Range("I66").Value = 2362.6428
$ws.Range("J66").Value = 7500.6
$ws.Range("K66").Value = 11813.214
$ws.Range("L66").Value = 37503
$ws.Range("M66").Value = -8381.214
$ws.Range("N66").Value = -44367
$ws.Range("H74").Value = 83925.375
$ws.Range("I74").Value = 57649.26
$ws.Range("K74").Value = 57649.26
$ws.Range("M74").Value = -56775.26
$ws.Range("H77").Value = 83925.375
$ws.Range("I77").Value = 57649.26
$ws.Range("K77").Value = 288246.3
$ws.Range("M77").Value = -283878.3
$ws.Range("H102").Value = 8335893
$ws.Range("I102").Value = 10418502
$ws.Range("K102").Value = 10418502
$ws.Range("M102").Value = -10416880
$ws.Range("H116").Value = 2828749.8
$ws.Range("I116").Value = 3770928.2
$ws.Range("K116").Value = 3770928.2
$ws.Range("M116").Value = -3768634.2
$ws.Range("H122").Value = 10501697
$ws.Range("J122").Value = 4175260.5
$ws.Range("L122").Value = 12525781.5
$ws.Range("N122").Value = -12530681.5
$ws.Range("H132").Value = 4527.171
$ws.Range("I132").Value = 4626.5
$ws.Range("J132").Value = 4256.273
$ws.Range("K132").Value = 13879.5
$ws.Range("L132").Value = 12768.819
$ws.Range("M132").Value = -11349.5
$ws.Range("N132").Value = -17828.819
$ws.Range("H136").Value = 7349.9287
$ws.Range("I136").Value = 7531.077
$ws.Range("K136").Value = 22593.231
$ws.Range("M136").Value = -20043.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2828749.8
$ws.Range("I3").Value = 3770928.2
$ws.Range("K3").Value = 3770928.2
$ws.Range("M3").Value = -3770814.2
$ws.Range("H64").Value = 1119.0769
$ws.Range("I64").Value = 475
$ws.Range("J64").Value = 1236.1818
$ws.Range("K64").Value = 475
$ws.Range("L64").Value = 1236.1818
$ws.Range("M64").Value = -250
$ws.Range("N64").Value = -1686.1818
$ws.Range("H67").Value = 1119.0769
$ws.Range("I67").Value = 475
$ws.Range("J67").Value = 1236.1818
$ws.Range("K67").Value = 475
$ws.Range("L67").Value = 1236.1818
$ws.Range("M67").Value = 305
$ws.Range("N67").Value = -2796.1818
$ws.Range("H86").Value = 10014057
$ws.Range("I86").Value = 14290224
$ws.Range("K86").Value = 14290224
$ws.Range("M86").Value = -14289101
$ws.Range("H89").Value = 10014057
$ws.Range("I89").Value = 14290224
$ws.Range("K89").Value = 71451120
$ws.Range("M89").Value = -71445504

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28980.732
$ws.Range("I31").Value = 9762.666999999999
$ws.Range("K31").Value = 9762.666999999999
$ws.Range("M31").Value = -9467.666999999999
$ws.Range("H34").Value = 28980.732
$ws.Range("I34").Value = 9762.666999999999
$ws.Range("K34").Value = 9762.666999999999
$ws.Range("M34").Value = -9560.666999999999
$ws.Range("H86").Value = 11534.488
$ws.Range("I86").Value = 10302.552
$ws.Range("K86").Value = 10302.552
$ws.Range("M86").Value = -9179.552
$ws.Range("H89").Value = 11534.488
$ws.Range("I89").Value = 10302.552
$ws.Range("K89").Value = 51512.75999999999
$ws.Range("M89").Value = -45896.75999999999
$ws.Range("H132").Value = 127632.93
$ws.Range("I132").Value = 84968
$ws.Range("J132").Value = 298292.66
$ws.Range("K132").Value = 254904
$ws.Range("L132").Value = 894877.98
$ws.Range("M132").Value = -252374
$ws.Range("N132").Value = -899937.98

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2242
$ws.Range("I132").Value = 1812.25
$ws.Range("J132").Value = 2815
$ws.Range("K132").Value = 16310.25
$ws.Range("L132").Value = 25335
$ws.Range("M132").Value = -13780.25
$ws.Range("N132").Value = -30395
$ws.Range("H139").Value = 1884.6666
$ws.Range("I139").Value = 1884.6666
$ws.Range("K139").Value = 5653.9998
$ws.Range("M139").Value = -513.9997999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33337352
$ws.Range("J70").Value = 4853.5
$ws.Range("L70").Value = 4853.5
$ws.Range("N70").Value = -5393.5
$ws.Range("H73").Value = 33337352
$ws.Range("J73").Value = 4853.5
$ws.Range("L73").Value = 4853.5
$ws.Range("N73").Value = -6725.5
$ws.Range("H122").Value = 300534.44
$ws.Range("I122").Value = 426694.38
$ws.Range("K122").Value = 1280083.14
$ws.Range("M122").Value = -1277633.14
$ws.Range("H132").Value = 5094.551
$ws.Range("I132").Value = 3926.52
$ws.Range("J132").Value = 8168.316
$ws.Range("K132").Value = 11779.56
$ws.Range("L132").Value = 24504.948
$ws.Range("M132").Value = -9249.559999999999
$ws.Range("N132").Value = -29564.948

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 236.90909
$ws.Range("I16").Value = 195.6
$ws.Range("K16").Value = 195.6
$ws.Range("M16").Value = -25.59999999999999
$ws.Range("H22").Value = 36926.88
$ws.Range("I22").Value = 47695.95
$ws.Range("K22").Value = 47695.95
$ws.Range("M22").Value = -47400.95
$ws.Range("H27").Value = 36926.88
$ws.Range("I27").Value = 47695.95
$ws.Range("K27").Value = 47695.95
$ws.Range("M27").Value = -47588.95
$ws.Range("H68").Value = 4354.636
$ws.Range("I68").Value = 4211.222
$ws.Range("K68").Value = 4211.222
$ws.Range("M68").Value = -3462.222
$ws.Range("H71").Value = 4354.636
$ws.Range("I71").Value = 4211.222
$ws.Range("K71").Value = 21056.11
$ws.Range("M71").Value = -17312.11
$ws.Range("H132").Value = 11769.023
$ws.Range("I132").Value = 12326.342
$ws.Range("J132").Value = 6474.5
$ws.Range("K132").Value = 36979.026
$ws.Range("L132").Value = 19423.5
$ws.Range("M132").Value = -34449.026
$ws.Range("N132").Value = -24483.5
$ws.Range("H136").Value = 32783.734
$ws.Range("I136").Value = 46675.82
$ws.Range("K136").Value = 140027.46
$ws.Range("M136").Value = -137477.46

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 7098
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H81").Value = 5955610.5
$ws.Range("J81").Value = 4576.25
$ws.Range("L81").Value = 9152.5
$ws.Range("N81").Value = -11274.5
$ws.Range("H84").Value = 5955610.5
$ws.Range("J84").Value = 4576.25
$ws.Range("L84").Value = 45762.5
$ws.Range("N84").Value = -56370.5
